$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("R3").Value = "Auflöseabschnitt (im Zielgleis)"
$ws.Range("S4").Value = "Bemessungslänge"
